# The commit swaps the presentation's applied theme from the "Integral" /
# "Red Violet" color scheme over to the stock "Office Theme" / "Office"
# color scheme (i.e. the user picked a different Design on the ribbon).
#
# ppt/theme/theme1.xml (the theme actually bound to the slide master / the
# slides you see) needs its 12 theme colors changed from the "Red Violet"
# palette to the standard "Office" palette. Drive that through the
# PowerPoint object model's ThemeColorScheme, which is exactly what the
# Design-gallery click updates under the hood.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target "Office" theme colors, in ThemeColorScheme order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgr = ($b * 65536) + ($g * 256) + $r

    $color = $tcs.Item($i)
    $color.RGB = $bgr
}
